# edit.ps1 - apply the "Add files via upload" commit:
#   1. Slide 2 (sldId 258), shape id 7 ("Segnaposto contenuto 6"): drop the
#      word "esperto" from the first sentence of the body placeholder.
#   2. Every cached "datetime1" field on the slide master and its 11
#      layouts is refreshed from 28/10/2020 to 03/11/2020, as PowerPoint
#      does whenever it resaves a deck and re-caches its auto date/time
#      placeholders.

$p = $ppt.ActivePresentation

# 1. Fix the body text on slide 2 -----------------------------------------
$slide = $p.Slides.Item(2)
$shape = $slide.Shapes.Item(2)
$shape.TextFrame.TextRange.Text = "Il sistema, partendo dalle preferenze musicali dell’utente, definite attraverso artisti e album ascoltati, consente di generare, studiando i generi di questi ultimi, classifiche di artisti simili a quelli inseriti dall’utente e playlists di canzoni, collegando l’account di Spotify."

# 2. Refresh the cached date placeholders ---------------------------------
$newDate = "03/11/2020"

$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

$layouts = $master.CustomLayouts
$layoutShapeIdx = @{1=4; 2=3; 3=4; 4=4; 5=6; 6=2; 7=1; 8=5; 9=4; 10=3; 11=7}
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $idx = $layoutShapeIdx[$li]
    $layout.Shapes.Item($idx).TextFrame.TextRange.Text = $newDate
}
